$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for Pomelo (Start Ruby / Primera).
# In the sheet this shows up as a new row inserted at row 83, which
# pushes the existing row 83 (and everything below it, through 190)
# down by one row, down to row 191.
$ws.Rows(83).Insert()

# Seed the newly-inserted (currently blank) row 83 with the metadata
# columns that stay the same across records (market, region, product,
# variety, quality, unit, origin, etc.) by copying them from the row
# immediately below, which now holds the record that used to be at 83.
$ws.Range("A84:T84").Copy()
$ws.Range("A83").PasteSpecial()

# Now overwrite just the columns that differ for this new record:
# date, volume, min/max/avg price and $/kg.
$ws.Range("D83").Value = 44546
$ws.Range("M83").Value = 200
$ws.Range("N83").Value = 11000
$ws.Range("O83").Value = 12000
$ws.Range("P83").Value = 11500
$ws.Range("S83").Value = 821
